# COREESG_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer (A10) from 2021-05-19 to 2021-05-20
#  - refresh the Weight (D) and Percent Change (E) figures for rows 2-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection afterwards.
$ws.Unprotect()

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."
# Re-fit row 10 so the wrapped disclaimer text doesn't leave a stale explicit
# row height behind.
$ws.Rows.Item(10).AutoFit()

$ws.Range("D2").Value = 0.2440787948252343
$ws.Range("E2").Value = 0.01842688960928163

$ws.Range("D3").Value = 0.5019317580930595
$ws.Range("E3").Value = 0.006102414433536874

$ws.Range("D4").Value = 0.09528326653460861
$ws.Range("E4").Value = 0.02244278081787332

$ws.Range("D5").Value = 0.1022716500471871
$ws.Range("E5").Value = 0.003567508232711347

$ws.Range("D6").Value = 0.05643453049991069
$ws.Range("E6").Value = 0.005983889528193176

$ws.Range("E7").Value = 0.01040158302949923

$ws.Protect()
